# [FIX] Use NISN for student import and handle schedule job
#
# The "NIS" header in the student-import template is renamed to "NISN"
# (column A), and the active selection is moved from E14 to A6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "NIS" -> "NISN"
$ws.Range("A1").Value = "NISN"

# Move/refresh the active selection
$ws.Range("A6").Select()
